# Applies the six text replacements described by the diff.
# Each replacement is performed by locating the unique old text with
# Find.Execute (no ReplaceWith argument, so no AutoFormat "smart quotes"
# substitution happens) and then assigning the new text straight onto the
# matched Range, which keeps punctuation (straight quotes, etc.) intact.

$d = $word.ActiveDocument

function Replace-UniqueText($oldText, $newText, $wholeWord) {
    $rng = $d.Content
    $found = $rng.Find.Execute($oldText, $true, $wholeWord, $false, $false, $false, $true, 1, $false)
    if ($found) {
        $rng.Text = $newText
    }
    return $found
}

# Contract number: 1234568 -> 2445987
Replace-UniqueText "1234568" "2445987" $false | Out-Null

# Day-of-month field "09" -> "01" (whole word match so the other
# occurrences of "09" that are part of larger numbers/dates, such as
# "09.01.2024" or "1027809210330", are left untouched).
Replace-UniqueText "09" "01" $true | Out-Null

# Month field: "февраля" -> "мая"
Replace-UniqueText "февраля" "мая" $false | Out-Null

# Signatory name (title block): "О. Ю." -> "О. Е."
Replace-UniqueText "генерального директора Котлярчука О. Ю." "генерального директора Котлярчука О. Е." $false | Out-Null

# Subject/description of the reviewed technical documentation
Replace-UniqueText "Рассмотрение технической документации ""Информация об остойчивости"" № 123-FU-TB-3455" "Рассмотрение технической документации ""Грузовая марка"" № 5234-234234-23 на т/х ""МУРМАН 1"" РС 091052" $false | Out-Null

# Signature line: "О. Ю. Котлярчук" -> "О. Е. Котлярчук"
Replace-UniqueText "О. Ю. Котлярчук" "О. Е. Котлярчук" $false | Out-Null

Write-Output "edits applied"
